# Edit: split the paragraph
#   "Bởi vì Thuộc tính của dimension rồi đến Members rồi lại đến Thuộc tính
#    của dimension (trong cây con) */"
# into three runs by inserting the phrase "dimension cùng cấp với " after
# "...rồi đến " and before "Members...", and relocate the document's
# "_GoBack" bookmark (collapsed / zero-length) to sit between the newly
# inserted run and "Members...".

$d = $word.ActiveDocument
$p = $d.Paragraphs(4)

# Locate "Members" inside this paragraph (scoped search so we don't touch
# any other occurrence elsewhere in the document).
$find = $p.Range.Duplicate
$found = $find.Find.Execute("Members", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'Members' in paragraph 4"
}

$splitPos = $find.Duplicate
$splitPos.Collapse(1)
$insertAt = $splitPos.Start

# Insert a placeholder run right before "Members" so we have a *non-collapsed*
# range to anchor a bookmark to (a bookmark added at a truly collapsed point
# gets pushed to the far side of whatever text is subsequently inserted
# there, instead of staying put / wrapping it).
$splitPos.InsertBefore("@@PLACEHOLDER@@")

$phRange = $d.Range($insertAt, $insertAt + ("@@PLACEHOLDER@@").Length)
$d.Bookmarks.Add("ZZ__tmp_split", $phRange) | Out-Null

# Replace the placeholder text with the real inserted phrase. Doing the
# replacement through the bookmark's own Range keeps the new text properly
# enclosed by the (temporary) bookmark, which in turn keeps it as a
# separate run from the text before it once the temporary bookmark is
# removed again.
$phRange2 = $d.Bookmarks("ZZ__tmp_split").Range
$phRange2.Text = "dimension cùng cấp với "

# Collapse to the end of the newly-inserted text (i.e. right before
# "Members...") and drop the real "_GoBack" bookmark there. Word only ever
# keeps a single "_GoBack" bookmark, so adding this one automatically
# removes it from wherever it used to be (the near-empty paragraph further
# down in the document).
$goBackPoint = $d.Bookmarks("ZZ__tmp_split").Range.Duplicate
$goBackPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# Clean up the temporary bookmark; the run split it was protecting remains.
$d.Bookmarks("ZZ__tmp_split").Delete()
